$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (49 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1954.4
$ws.Range("J17").Value = 1954.4
$ws.Range("L17").Value = 5863.200000000001
$ws.Range("N17").Value = -6199.200000000001
$ws.Range("H29").Value = 4987.4
$ws.Range("I29").Value = 437
$ws.Range("J29").Value = 6125
$ws.Range("K29").Value = 1311
$ws.Range("L29").Value = 18375
$ws.Range("M29").Value = -1030
$ws.Range("N29").Value = -18937
$ws.Range("H70").Value = 14758.546
$ws.Range("I70").Value = 2249
$ws.Range("J70").Value = 19449.625
$ws.Range("K70").Value = 6747
$ws.Range("L70").Value = 58348.875
$ws.Range("M70").Value = -6477
$ws.Range("N70").Value = -58888.875
$ws.Range("H73").Value = 14758.546
$ws.Range("I73").Value = 2249
$ws.Range("J73").Value = 19449.625
$ws.Range("K73").Value = 6747
$ws.Range("L73").Value = 58348.875
$ws.Range("M73").Value = -5811
$ws.Range("N73").Value = -60220.875
$ws.Range("H87").Value = 69994.55499999999
$ws.Range("J87").Value = 87658.5
$ws.Range("L87").Value = 87658.5
$ws.Range("N87").Value = -90154.5
$ws.Range("H90").Value = 69994.55499999999
$ws.Range("J90").Value = 87658.5
$ws.Range("L90").Value = 262975.5
$ws.Range("N90").Value = -275455.5
$ws.Range("H112").Value = 1426.9667
$ws.Range("J112").Value = 1491.6296
$ws.Range("L112").Value = 4474.8888
$ws.Range("N112").Value = -6690.8888
$ws.Range("H133").Value = 80775
$ws.Range("J133").Value = 80775
$ws.Range("L133").Value = 80775
$ws.Range("N133").Value = -90895
$ws.Range("H136").Value = 94998
$ws.Range("J136").Value = 94998
$ws.Range("L136").Value = 94998
$ws.Range("N136").Value = -105198
$ws.Range("H138").Value = 2543.173
$ws.Range("J138").Value = 3712.8262
$ws.Range("L138").Value = 11138.4786
$ws.Range("N138").Value = -21418.4786

# --- Sheet: ARM (28 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5960.5884
$ws.Range("I61").Value = 6555.9165
$ws.Range("K61").Value = 6555.9165
$ws.Range("M61").Value = -6343.9165
$ws.Range("H64").Value = 71199.39999999999
$ws.Range("J64").Value = 72666
$ws.Range("L64").Value = 72666
$ws.Range("N64").Value = -73162
$ws.Range("H67").Value = 71199.39999999999
$ws.Range("J67").Value = 72666
$ws.Range("L67").Value = 72666
$ws.Range("N67").Value = -74382
$ws.Range("H110").Value = 803.75
$ws.Range("I110").Value = 825.13336
$ws.Range("K110").Value = 825.13336
$ws.Range("M110").Value = 1219.86664
$ws.Range("H132").Value = 3169.158
$ws.Range("I132").Value = 2296.26
$ws.Range("K132").Value = 6888.780000000001
$ws.Range("M132").Value = -4358.780000000001
$ws.Range("H136").Value = 5960.5884
$ws.Range("I136").Value = 6555.9165
$ws.Range("K136").Value = 19667.7495
$ws.Range("M136").Value = -17117.7495
$ws.Range("H139").Value = 84999.836
$ws.Range("J139").Value = 84999.836
$ws.Range("L139").Value = 84999.836
$ws.Range("N139").Value = -95279.836

# --- Sheet: BSM (20 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 74090.5
$ws.Range("J62").Value = 74090.5
$ws.Range("L62").Value = 74090.5
$ws.Range("N62").Value = -75462.5
$ws.Range("H65").Value = 74090.5
$ws.Range("J65").Value = 74090.5
$ws.Range("L65").Value = 222271.5
$ws.Range("N65").Value = -229135.5
$ws.Range("H92").Value = 67542.14
$ws.Range("J92").Value = 67542.14
$ws.Range("L92").Value = 67542.14
$ws.Range("N92").Value = -72534.14
$ws.Range("H107").Value = 4169.9
$ws.Range("I107").Value = 3673.3333
$ws.Range("K107").Value = 3673.3333
$ws.Range("M107").Value = -1753.3333
$ws.Range("H140").Value = 166913.4
$ws.Range("J140").Value = 166913.4
$ws.Range("L140").Value = 166913.4
$ws.Range("N140").Value = -177273.4

# --- Sheet: CRP (41 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4789.1304
$ws.Range("I31").Value = 3349.05
$ws.Range("J31").Value = 5896.885
$ws.Range("K31").Value = 3349.05
$ws.Range("L31").Value = 5896.885
$ws.Range("M31").Value = -3054.05
$ws.Range("N31").Value = -6486.885
$ws.Range("H34").Value = 4789.1304
$ws.Range("I34").Value = 3349.05
$ws.Range("J34").Value = 5896.885
$ws.Range("K34").Value = 3349.05
$ws.Range("L34").Value = 5896.885
$ws.Range("M34").Value = -3147.05
$ws.Range("N34").Value = -6300.885
$ws.Range("H51").Value = 39285.285
$ws.Range("J51").Value = 39285.285
$ws.Range("L51").Value = 39285.285
$ws.Range("N51").Value = -40757.285
$ws.Range("H61").Value = 39285.285
$ws.Range("J61").Value = 39285.285
$ws.Range("L61").Value = 39285.285
$ws.Range("N61").Value = -39981.285
$ws.Range("H68").Value = 67184.92999999999
$ws.Range("J68").Value = 68583.766
$ws.Range("L68").Value = 68583.766
$ws.Range("N68").Value = -70081.766
$ws.Range("H71").Value = 67184.92999999999
$ws.Range("J71").Value = 68583.766
$ws.Range("L71").Value = 205751.298
$ws.Range("N71").Value = -213239.298
$ws.Range("H122").Value = 3024.8235
$ws.Range("I122").Value = 1324
$ws.Range("J122").Value = 4215.4
$ws.Range("K122").Value = 3972
$ws.Range("L122").Value = 12646.2
$ws.Range("M122").Value = -1522
$ws.Range("N122").Value = -17546.2
$ws.Range("H134").Value = 1966.9565
$ws.Range("I134").Value = 1933.7368
$ws.Range("K134").Value = 5801.2104
$ws.Range("M134").Value = -3266.2104

# --- Sheet: CUL (18 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 1100
$ws.Range("J21").Value = 1100
$ws.Range("L21").Value = 3300
$ws.Range("N21").Value = -3646
$ws.Range("H75").Value = 4926.7334
$ws.Range("I75").Value = 998.5
$ws.Range("J75").Value = 5531.077
$ws.Range("K75").Value = 2995.5
$ws.Range("L75").Value = 16593.231
$ws.Range("M75").Value = -1997.5
$ws.Range("N75").Value = -18589.231
$ws.Range("H78").Value = 4926.7334
$ws.Range("I78").Value = 998.5
$ws.Range("J78").Value = 5531.077
$ws.Range("K78").Value = 8986.5
$ws.Range("L78").Value = 49779.693
$ws.Range("M78").Value = -3994.5
$ws.Range("N78").Value = -59763.693

# --- Sheet: GSM (19 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1857279.5
$ws.Range("I14").Value = 1537694.2
$ws.Range("J14").Value = 2496450
$ws.Range("K14").Value = 1537694.2
$ws.Range("L14").Value = 2496450
$ws.Range("M14").Value = -1537526.2
$ws.Range("N14").Value = -2496786
$ws.Range("H20").Value = 41560.2
$ws.Range("J20").Value = 41560.2
$ws.Range("L20").Value = 41560.2
$ws.Range("N20").Value = -42050.2
$ws.Range("H97").Value = 1720.7188
$ws.Range("I97").Value = 671.1818
$ws.Range("K97").Value = 671.1818
$ws.Range("M97").Value = -175.1818
$ws.Range("H122").Value = 4006.1538
$ws.Range("I122").Value = 2924.25
$ws.Range("K122").Value = 8772.75
$ws.Range("M122").Value = -6322.75

# --- Sheet: LTW (8 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 7312.381
$ws.Range("I100").Value = 2313.3333
$ws.Range("K100").Value = 2313.3333
$ws.Range("M100").Value = -1772.3333
$ws.Range("H116").Value = 80543.8
$ws.Range("J116").Value = 80543.8
$ws.Range("L116").Value = 80543.8
$ws.Range("N116").Value = -89721.8

# --- Sheet: WVR (4 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1215.3334
$ws.Range("J96").Value = 918.75
$ws.Range("L96").Value = 918.75
$ws.Range("N96").Value = -3664.75

Write-Output "Applied 187 cell updates across 8 sheets"
